# EPBDS-11051 Improve performance in OpenLClassLoader usage
#
# The "Rules" sheet contained a leftover example/test block (rows 18-21,
# columns B:C) describing a "mapSpr" spreadsheet rule - its title, the
# "Steps"/"Values" header and two step rows. This leftover example data
# is no longer needed, so its contents are cleared out while the cells
# themselves and their existing formatting (the quote-prefixed style
# applied to C20:C21) are left in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover "mapSpr" example text from B18:C21 without
# touching the cell formatting/styles already applied to that range.
$ws.Range("B18:C21").ClearContents()

# The user ends up with A16 selected instead of the old C20 selection.
$ws.Range("A16").Select()

# Remember the on-screen window position at save time.
$win = $wb.Windows.Item(1)
$win.Left = 12940
$win.Top = 3130

$wb.Save()
